$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for B2:I21 (regenerated after adding a fixed random seed)
$data = @(
    @(0, 0.02, 0.04, -0.04, 0, 0.03, 0.07000000000000001, -0.07000000000000001),
    @(0.01, 0.03, 0.09, -0.05, 0.01, 0.05, 0.13, -0.11),
    @(0, 0.02, 0.05, -0.05, 0, 0.03, 0.07000000000000001, -0.08),
    @(0, 0.02, 0.04, -0.03, 0, 0.02, 0.06, -0.06),
    @(0, 0.03, 0.07000000000000001, -0.07000000000000001, 0, 0.05, 0.13, -0.13),
    @(0, 0.02, 0.06, -0.07000000000000001, -0.01, 0.04, 0.11, -0.1),
    @(0, 0.03, 0.09, -0.06, 0, 0.04, 0.1, -0.11),
    @(0, 0.02, 0.07000000000000001, -0.07000000000000001, 0, 0.04, 0.11, -0.13),
    @(0, 0.02, 0.04, -0.05, 0, 0.03, 0.08, -0.07000000000000001),
    @(0, 0.02, 0.05, -0.05, 0, 0.03, 0.08, -0.1),
    @(0, 0.02, 0.05, -0.04, 0, 0.03, 0.07000000000000001, -0.08),
    @(0, 0.02, 0.04, -0.05, -0.01, 0.03, 0.08, -0.12),
    @(0, 0.02, 0.04, -0.04, 0, 0.03, 0.08, -0.1),
    @(0, 0.02, 0.04, -0.04, 0, 0.03, 0.07000000000000001, -0.09),
    @(0, 0.01, 0.04, -0.04, 0, 0.02, 0.08, -0.05),
    @(0, 0.02, 0.06, -0.05, 0, 0.03, 0.07000000000000001, -0.1),
    @(0, 0.04, 0.08, -0.09, -0.01, 0.05, 0.14, -0.13),
    @(0, 0.02, 0.06, -0.04, 0, 0.03, 0.07000000000000001, -0.08),
    @(0, 0.03, 0.1, -0.06, 0, 0.05, 0.12, -0.13),
    @(0, 0.03, 0.07000000000000001, -0.07000000000000001, 0, 0.05, 0.14, -0.11)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = 2 + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $c = 2 + $j
        $ws.Cells.Item($r, $c).Value = $rowVals[$j]
    }
}
